# Daily II - poprawka
# Poprawiona wersja dokumentu z opisem pracy przy FEAT
#
# Kamil Sajdak's row (row 4) on sheet "Arkusz2" is updated:
#  - "Co zrobiłam/em?" (col B) gets a new description of the completed work
#  - "Jakie mam/miałem problemy?" (col D) gets a new description (was "-")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz2")

$ws.Range("B4").Value = "Naniesienie poprawek do dokumentu z cechami systemu FEAT, dodanie kilku nowych cech."
$ws.Range("D4").Value = "Przeanalizowanie dokumentów w poszukiwaniu nowych cech systemu."

# Row grew taller to fit the new wrapped text
$ws.Rows.Item(4).RowHeight = 39.6

# Leave the cursor/selection on the cell that was last edited
$ws.Range("D4").Select() | Out-Null
